# Update countries & provincias Spain
# - Refresh the "last updated" timestamp.
# - Re-sort several country rows (name + stats move together) to reflect
#   the new case-count ordering (Indonesia, Estonia/Bielorrusia,
#   Islas Caimanes, Laos, Suazilandia, the Turcas/Somalia/San Vicente/
#   Malaui block, and Islas Virgenes Britanicas all shift position).
# - Refresh a handful of per-country statistic values (rows 18, 32, 111)
#   that were updated independently of any re-sort.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Datos actualizados a 8 de Abril de 2020 a las 10:52"

$ws.Range("D18").Value = 4512
$ws.Range("E18").Value = 7936
$ws.Range("F18").Value = 267
$ws.Range("G18").Value = 30
$ws.Range("H18").Value = 273

$ws.Range("E32").Value = 3750
$ws.Range("G32").Value = 10
$ws.Range("H32").Value = 207

$ws.Range("A39").Value = "Indonesia"
$ws.Range("B39").Value = 2956
$ws.Range("C39").Value = 218
$ws.Range("D39").Value = 222
$ws.Range("E39").Value = 2494
$ws.Range("F39").Value = 0
$ws.Range("G39").Value = 19
$ws.Range("H39").Value = 240

$ws.Range("A40").Value = "Peru"
$ws.Range("B40").Value = 2954
$ws.Range("D40").Value = 1301
$ws.Range("E40").Value = 1546
$ws.Range("F40").Value = 109
$ws.Range("H40").Value = 107

$ws.Range("A41").Value = "Arabia Saudita"
$ws.Range("B41").Value = 2795
$ws.Range("C41").Value = 0
$ws.Range("D41").Value = 615
$ws.Range("E41").Value = 2139
$ws.Range("F41").Value = 41
$ws.Range("G41").Value = 0
$ws.Range("H41").Value = 41

$ws.Range("A42").Value = "Mexico"
$ws.Range("B42").Value = 2785
$ws.Range("C42").Value = 346
$ws.Range("D42").Value = 633
$ws.Range("E42").Value = 2011
$ws.Range("F42").Value = 89
$ws.Range("G42").Value = 16
$ws.Range("H42").Value = 141

$ws.Range("A61").Value = "Estonia"
$ws.Range("B61").Value = 1185
$ws.Range("C61").Value = 36
$ws.Range("D61").Value = 72
$ws.Range("E61").Value = 1089
$ws.Range("F61").Value = 11
$ws.Range("G61").Value = 3
$ws.Range("H61").Value = 24

$ws.Range("A62").Value = "Marruecos"
$ws.Range("B62").Value = 1184
$ws.Range("D62").Value = 93
$ws.Range("E62").Value = 1001
$ws.Range("F62").Value = 1
$ws.Range("H62").Value = 90

$ws.Range("A64").Value = "Bielorrusia"
$ws.Range("B64").Value = 1066
$ws.Range("C64").Value = 205
$ws.Range("D64").Value = 77
$ws.Range("E64").Value = 976
$ws.Range("F64").Value = 33
$ws.Range("H64").Value = 13

$ws.Range("A65").Value = "Eslovenia"
$ws.Range("B65").Value = 1059
$ws.Range("D65").Value = 102
$ws.Range("E65").Value = 921
$ws.Range("F65").Value = 30
$ws.Range("G65").Value = 0
$ws.Range("H65").Value = 36

$ws.Range("A66").Value = "Moldavia"
$ws.Range("B66").Value = 1056
$ws.Range("D66").Value = 40
$ws.Range("E66").Value = 992
$ws.Range("F66").Value = 80
$ws.Range("G66").Value = 2
$ws.Range("H66").Value = 24

$ws.Range("A67").Value = "Hong Kong"
$ws.Range("B67").Value = 961
$ws.Range("C67").Value = 25
$ws.Range("D67").Value = 264
$ws.Range("E67").Value = 693
$ws.Range("F67").Value = 14
$ws.Range("H67").Value = 4

$ws.Range("A68").Value = "Lituania"
$ws.Range("B68").Value = 912
$ws.Range("C68").Value = 32
$ws.Range("D68").Value = 8
$ws.Range("E68").Value = 889
$ws.Range("F68").Value = 11
$ws.Range("G68").Value = 0
$ws.Range("H68").Value = 15

$ws.Range("A69").Value = "Hungria"
$ws.Range("B69").Value = 895
$ws.Range("C69").Value = 78
$ws.Range("D69").Value = 94
$ws.Range("E69").Value = 743
$ws.Range("F69").Value = 17
$ws.Range("G69").Value = 11
$ws.Range("H69").Value = 58

$ws.Range("A70").Value = "Armenia"
$ws.Range("B70").Value = 881
$ws.Range("C70").Value = 28
$ws.Range("D70").Value = 114
$ws.Range("E70").Value = 758
$ws.Range("F70").Value = 30
$ws.Range("G70").Value = 1
$ws.Range("H70").Value = 9

$ws.Range("D111").Value = 48
$ws.Range("E111").Value = 157

$ws.Range("A144").Value = "Islas Caimanes"
$ws.Range("D144").Value = 6
$ws.Range("H144").Value = 1

$ws.Range("A145").Value = "Congo"
$ws.Range("D145").Value = 2
$ws.Range("H145").Value = 5

$ws.Range("A173").Value = "Laos"
$ws.Range("C173").Value = 1
$ws.Range("D173").Value = 0
$ws.Range("E173").Value = 15

$ws.Range("A174").Value = "Dominica"
$ws.Range("B174").Value = 15
$ws.Range("D174").Value = 1

$ws.Range("A186").Value = "Suazilandia"
$ws.Range("D186").Value = 4
$ws.Range("H186").Value = 0

$ws.Range("A187").Value = "Surinam"
$ws.Range("D187").Value = 3
$ws.Range("H187").Value = 1

$ws.Range("A191").Value = "Islas Turcas y Caicos"
$ws.Range("F191").Value = 0

$ws.Range("A192").Value = "Somalia"
$ws.Range("D192").Value = 1
$ws.Range("H192").Value = 0

$ws.Range("A194").Value = "Malaui"
$ws.Range("D194").Value = 0
$ws.Range("F194").Value = 1
$ws.Range("H194").Value = 1

$ws.Range("A208").Value = "Islas Virgenes Britanicas"

$ws.Range("A209").Value = "Anguila"
